# Update control flow ppt slides
# - lower_snake_case the Write/ReadLine/WriteLine API calls used in the
#   on-slide "code" textboxes (Write -> write, ReadLine -> read_line,
#   WriteLine -> write_line)
# - change the language the demo chats about from C# to C++ in the code
#   textboxes (slides 1-4) and the standalone "language" labels, while
#   slides 5-6 (which previously said C++) get reverted to C#
# - keep the "Good choice, ..." sentence consistent with the new language
#   ("a fine language" -> "a great language")

function Replace-SubText($TextRange, $Old, $New) {
    $full = $TextRange.Text
    $idx = $full.IndexOf($Old)
    if ($idx -lt 0) {
        return $false
    }
    $start = $idx + 1
    $len = $Old.Length
    $sub = $TextRange.Characters($start, $len)
    $sub.Text = $New
    return $true
}

function Find-ShapeByName($container, $name) {
    try {
        $shp = $container.Shapes.Item($name)
        return $shp
    } catch {
    }
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $top = $container.Shapes.Item($i)
        if ($top.Type -eq 6) {
            try {
                $found = $top.GroupItems.Item($name)
                return $found
            } catch {
            }
        }
    }
    return $null
}

function Update-CodeBox($slide, $pairs) {
    $shp = Find-ShapeByName $slide "TextBox 15"
    if ($shp -eq $null) {
        Write-Host "code textbox not found"
        return
    }
    $tr = $shp.TextFrame.TextRange
    foreach ($pair in $pairs) {
        $ok = Replace-SubText $tr $pair[0] $pair[1]
        if (-not $ok) {
            Write-Host "NOT FOUND in code box: $($pair[0])"
        }
    }
}

function Update-SimpleBox($slide, $name, $old, $new) {
    $shp = Find-ShapeByName $slide $name
    if ($shp -eq $null) {
        Write-Host "$name not found"
        return
    }
    $tr = $shp.TextFrame.TextRange
    if ($tr.Text -ne $old) {
        Write-Host "$name text mismatch, got [$($tr.Text)] expected [$old]"
        return
    }
    $ok = Replace-SubText $tr $old $new
    if (-not $ok) {
        Write-Host "failed to replace text in $name"
    }
}

$p = $ppt.ActivePresentation

# Code block pairs shared by slides 1-5 (single combined "WriteLine(...)" run)
$codePairsCombined = @(
    ,@('Write("What language do you use? ");', 'write("What language do you use? ");')
    ,@('ReadLine', 'read_line')
    ,@('if (language == "C#")', 'if (language == "C++")')
    ,@('    WriteLine("Good choice, C# is a fine language.");', '    write_line("Good choice, C++ is a great language.");')
    ,@('WriteLine("Great chat!");', 'write_line("Great chat!");')
)

# Code block pairs for slide 6, where "WriteLine" and the "(...)" call are
# split across two separate runs
$codePairsSplit = @(
    ,@('Write("What language do you use? ");', 'write("What language do you use? ");')
    ,@('ReadLine', 'read_line')
    ,@('if (language == "C#")', 'if (language == "C++")')
    ,@('    WriteLine', '    write_line')
    ,@('("Good choice, C# is a fine language.");', '("Good choice, C++ is a great language.");')
    ,@('WriteLine("Great chat!");', 'write_line("Great chat!");')
)

# Slide 1
$s1 = $p.Slides.Item(1)
Update-CodeBox $s1 $codePairsCombined
Update-SimpleBox $s1 "TextBox 61" "C#" "C++"
Update-SimpleBox $s1 "TextBox 64" "C#" "C++"

# Slide 2
$s2 = $p.Slides.Item(2)
Update-CodeBox $s2 $codePairsCombined
Update-SimpleBox $s2 "TextBox 61" "C#" "C++"
Update-SimpleBox $s2 "TextBox 64" "C#" "C++"

# Slide 3
$s3 = $p.Slides.Item(3)
Update-CodeBox $s3 $codePairsCombined
Update-SimpleBox $s3 "TextBox 61" "C#" "C++"
Update-SimpleBox $s3 "TextBox 27" "C#" "C++"
Update-SimpleBox $s3 "TextBox 29" "Good choice, C# is a fine language." "Good choice, C++ is a great language."

# Slide 4
$s4 = $p.Slides.Item(4)
Update-CodeBox $s4 $codePairsCombined
Update-SimpleBox $s4 "TextBox 61" "C#" "C++"
Update-SimpleBox $s4 "TextBox 64" "C#" "C++"
Update-SimpleBox $s4 "TextBox 26" "Good choice, C# is a fine language." "Good choice, C++ is a great language."

# Slide 5
$s5 = $p.Slides.Item(5)
Update-CodeBox $s5 $codePairsCombined
Update-SimpleBox $s5 "TextBox 61" "C++" "C#"
Update-SimpleBox $s5 "TextBox 64" "C++" "C#"

# Slide 6
$s6 = $p.Slides.Item(6)
Update-CodeBox $s6 $codePairsSplit
Update-SimpleBox $s6 "TextBox 61" "C++" "C#"
Update-SimpleBox $s6 "TextBox 64" "C++" "C#"

Write-Host "done"
